# Add a new worksheet "Feuil2" after the existing "Feuil1"
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Feuil2"

# Put "Content" in E13, centered with a yellow fill
$cell = $newSheet.Range("E13")
$cell.Value = "Content"
$cell.HorizontalAlignment = -4108  # xlCenter
$cell.Interior.Color = 65535       # yellow (RGB FFFF00)

# Select F13 on the new sheet and make it the active sheet/tab
$newSheet.Range("F13").Select()
$newSheet.Activate()
